# Corrige les numéros de groupes : 030502XXXX -> 040311XXXX pour les paroisses
# (colonne C, lignes 2 a 16), puis repositionne la cellule active.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newGroupNumbers = [ordered]@{
    2  = 403110100
    3  = 403110200
    4  = 403110300
    5  = 403110400
    6  = 403110500
    7  = 403110600
    8  = 403110700
    9  = 403110800
    10 = 403110900
    11 = 403111000
    12 = 403111100
    13 = 403111200
    14 = 403111300
    15 = 403111400
    16 = 403111500
}

foreach ($row in $newGroupNumbers.Keys) {
    $cell = $ws.Cells.Item($row, 3)   # colonne C = IdxSG
    $cell.Value = $newGroupNumbers[$row]
    $cell.ClearFormats()
}

$ws.Range("K9").Select() | Out-Null
